$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated simulated error-structure values (triangle data, rows 2-25, cols A-L)
$data = @{
    2 = @(0.000168894302993827, 0.00195936044226156, 0.00466532218596389, 0.00702293731362239, 0.00856029066093903, 0.0104675215553273, 0.0125351872339614, 0.0155007238767398, 0.0193428339687708, 0.0229369923545958, 0.0263236012926186, 0.0293301372184897)
    3 = @(0.00308963739616441, 0.00609154701822014, 0.00854808528826605, 0.0100646642824376, 0.0120039410870029, 0.0141125299186969, 0.0171675226292971, 0.0211281524767944, 0.0248086275856852, 0.0282582129813217, 0.0313026578880072, 0.0340365679057244)
    4 = @(0.00213448800127303, 0.0030990025548629, 0.00327002746927913, 0.0042138669481561, 0.00554978676786286, 0.00822851088280566, 0.0119977363870131, 0.0153166515318061, 0.01838200447226, 0.0210269753861671, 0.0233804772423839, 0.0253842126988021)
    5 = @(0.000258527145586385, 0.00148837240581924, 0.00128086934320029, 0.00125100359818044, 0.00337198997451308, 0.00689281851646095, 0.00976195838470019, 0.0124106585795832, 0.0146567137602003, 0.0166510432888408, 0.0183210385429129, 0.0195049661484962)
    6 = @(0.00159132413079721, 0.00112525298711098, 0.00147709191948516, 0.00423801988727446, 0.00816522500253734, 0.0112637750194319, 0.0140654394694266, 0.0164075356394216, 0.0184660147201351, 0.0201743503069413, 0.0213709807956511, 0.0224975781967052)
    7 = @(0.00304906856854359, 0.00492141360594618, 0.00865927629341475, 0.0133667499452153, 0.0170705019881673, 0.0203377025321152, 0.0230445266622085, 0.0253968449994093, 0.0273435836315193, 0.0287297028395091, 0.0300186637727774, 0.0308782520494861)
    8 = @(0.000417279035482387, 0.00367999844192774, 0.00766544102921716, 0.0103530508547945, 0.0127225236261034, 0.0146124914599678, 0.016250573445635, 0.017561682499703, 0.0183660523077062, 0.0191522427453285, 0.0195589620975683, 0.0201381370330262)
    9 = @(0.00438642506966147, 0.00844045524149785, 0.0108547130070656, 0.0129950427360719, 0.0146660950879457, 0.016115364248265, 0.0172532607322062, 0.0178867782360055, 0.0185315227726725, 0.0188030303073172, 0.0192741975892211, 0.0201047443847868)
    10 = @(0.0027081187188811, 0.0027701110242512, 0.00307483527744122, 0.00320449802154213, 0.00334720435921186, 0.00336741169686765, 0.00315514708255138, 0.00300849071867399, 0.0028561187380671, 0.00272232357294147, 0.00293866351847164, 0.00450401998932132)
    11 = @(0.00235995840858028, 0.00305974291756069, 0.00396550755008332, 0.00468387599367417, 0.00546119031945916, 0.0066790525817217, 0.00749663764636224, 0.00867036677424353, 0.00925728781074532, 0.00929944117457032, 0.00889473227467554, 0.00876201195435931)
    12 = @(0.000923452485141141, 0.000968775642536291, 0.00114246331700943, 0.00114163120364033, 0.00111028721916654, 0.0010372367991315, 0.00152546054039194, 0.00149425609393288, 0.00163652695965731, 0.00356473942111562, 0.0068908381443954, 0.00960724061830813)
    13 = @(0.000758240260190668, 0.00090910549686826, 0.00131885794036515, 0.00250189028000069, 0.00297332189865942, 0.00402042997220463, 0.004285182949799, 0.00407355511882347, 0.00424832356907521, 0.00623976456325654, 0.00831772323724421, 0.0106731793709106)
    14 = @(0.00041408374922991, 0.000307132649023268, 0.00116616086340742, 0.00126670415076534, 0.00216389064298635, 0.00217648812058101, 0.00209515682299753, 0.00378842695778932, 0.00724522259419792, 0.0100590997536929, 0.0129063472502608, 0.0152571008936402)
    15 = @(0.000663702982715897, 0.00222608838266986, 0.00246462814532986, 0.00355172178717812, 0.00363588517422753, 0.00332176991580359, 0.00411177902944942, 0.00710013790623454, 0.00971670755501868, 0.0124568487239405, 0.0147147955316703)
    16 = @(0.00180708701790078, 0.00150161752170172, 0.00236860372900485, 0.00217859102174147, 0.00222769112309553, 0.00462247291084334, 0.00868341370548235, 0.0118412608878156, 0.0149526772372857, 0.0174692021575931)
    17 = @(0.00234572251214971, 0.00189452237371027, 0.00298121150786953, 0.00525848759390445, 0.00916976750898478, 0.0141622525103037, 0.0180448528607005, 0.0217121181467852, 0.0246703067965552)
    18 = @(0.00319261839024598, 0.00259735090163035, 0.00229950283528988, 0.00485547378713808, 0.00928248049945745, 0.0124819581402917, 0.0155943249543359, 0.0180206614633209)
    19 = @(0.00427110743552461, 0.00768053237826407, 0.0127297348157325, 0.0186826203349196, 0.023039558884586, 0.0270501338738557, 0.0301984090715587)
    20 = @(0.0018749818068402, 0.00576251018138595, 0.0105610116656644, 0.0134249180764707, 0.0161757619225038, 0.018167288901964)
    21 = @(0.00437133199853168, 0.00865512169017708, 0.0105655760518433, 0.0125913997068072, 0.0139168119951123)
    22 = @(0.00311305415764807, 0.00262315077351242, 0.00278000762190523, 0.00259431520233616)
    23 = @(0.0039463621383149, 0.00466214897469657, 0.00614717109217871)
    24 = @(0.00228738196520517, 0.00225547661778539)
    25 = @(0.00211902574229317)
}

foreach ($r in $data.Keys) {
    $values = $data[$r]
    $n = $values.Count
    $arr = New-Object "object[,]" 1,$n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $lastCol = [char](64 + $n)
    $rng = $ws.Range("A$r" + ":" + "$lastCol$r")
    $rng.Value = $arr
}
